$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Shared-string text updates (Volume/Number header, report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/29/2024  Through  8/4/2024"

# --- Row 15 (Murder, 1st Precinct): numeric cells C15/G15/H15 become
#     text placeholders ("0" / "***.*") matching columns D15/E15 which
#     already hold those shared strings. Set value as text (leading
#     apostrophe) then copy the number format from a sibling "text"
#     cell so the cell keeps using the shared-string text style (s=14)
#     instead of Excel's auto-generated "@"/custom text style. ---
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null

$ws.Range("G15").Value = "'0"
$ws.Range("D15").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null

$ws.Range("H15").Value = "'***.*"
$ws.Range("E15").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null

# --- Bulk numeric updates for rows 16-26, 28 (new weekly crime counts) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -61.111111111111
$ws.Range("I16").Value = 73
$ws.Range("J16").Value = 85
$ws.Range("K16").Value = -14.117647058823
$ws.Range("L16").Value = -12.048192771084
$ws.Range("M16").Value = 73.809523809523
$ws.Range("N16").Value = -85.341365461847
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 12
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 129
$ws.Range("J17").Value = 81
$ws.Range("K17").Value = 59.259259259259
$ws.Range("L17").Value = 53.571428571428
$ws.Range("M17").Value = 230.769230769231
$ws.Range("N17").Value = 19.444444444444
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 139
$ws.Range("J18").Value = 111
$ws.Range("K18").Value = 25.225225225225
$ws.Range("L18").Value = -18.235294117647
$ws.Range("M18").Value = 33.653846153846
$ws.Range("N18").Value = -72.144288577154
$ws.Range("C19").Value = 16
$ws.Range("E19").Value = -48.387096774193
$ws.Range("F19").Value = 81
$ws.Range("G19").Value = 123
$ws.Range("H19").Value = -34.146341463414
$ws.Range("I19").Value = 651
$ws.Range("J19").Value = 721
$ws.Range("K19").Value = -9.708737864077
$ws.Range("L19").Value = -9.331476323119
$ws.Range("M19").Value = 4.326923076923
$ws.Range("N19").Value = -71.296296296296
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("I20").Value = 27
$ws.Range("J20").Value = 39
$ws.Range("K20").Value = -30.769230769230
$ws.Range("L20").Value = -34.146341463414
$ws.Range("M20").Value = 22.727272727272
$ws.Range("N20").Value = -94.953271028037
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -25.641025641025
$ws.Range("F21").Value = 119
$ws.Range("G21").Value = 171
$ws.Range("H21").Value = -30.409356725146
$ws.Range("I21").Value = 1033
$ws.Range("J21").Value = 1047
$ws.Range("K21").Value = -1.337153772683
$ws.Range("L21").Value = -6.936936936936
$ws.Range("M21").Value = 23.860911270983
$ws.Range("N21").Value = -73.681528662420
$ws.Range("C22").Value = 2
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -30
$ws.Range("I22").Value = 53
$ws.Range("J22").Value = 52
$ws.Range("K22").Value = 1.923076923076
$ws.Range("L22").Value = -15.873015873015
$ws.Range("M22").Value = 35.897435897435
$ws.Range("C24").Value = 83
$ws.Range("D24").Value = 97
$ws.Range("E24").Value = -14.432989690721
$ws.Range("F24").Value = 296
$ws.Range("G24").Value = 323
$ws.Range("H24").Value = -8.359133126934
$ws.Range("I24").Value = 2299
$ws.Range("J24").Value = 2280
$ws.Range("K24").Value = 0.833333333333
$ws.Range("L24").Value = -1.372801372801
$ws.Range("M24").Value = 122.771317829457
$ws.Range("C25").Value = 81
$ws.Range("D25").Value = 96
$ws.Range("E25").Value = -15.625
$ws.Range("F25").Value = 269
$ws.Range("G25").Value = 323
$ws.Range("H25").Value = -16.71826625387
$ws.Range("I25").Value = 2237
$ws.Range("J25").Value = 2308
$ws.Range("K25").Value = -3.076256499133
$ws.Range("L25").Value = -5.010615711252
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = -6.451612903225
$ws.Range("I26").Value = 264
$ws.Range("J26").Value = 219
$ws.Range("K26").Value = 20.547945205479
$ws.Range("L26").Value = 25.714285714285
$ws.Range("M26").Value = 79.591836734693
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 12
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 68
$ws.Range("J28").Value = 53
$ws.Range("K28").Value = 28.301886792452
$ws.Range("L28").Value = 4.615384615384

# --- Row 27 (UCR Rape*, 1st Precinct): same text-placeholder treatment
#     as row 15 above for C27/G27/H27. ---
$ws.Range("C27").Value = "'0"
$ws.Range("D27").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null

$ws.Range("G27").Value = "'0"
$ws.Range("D27").Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4122) | Out-Null

$ws.Range("H27").Value = "'***.*"
$ws.Range("E27").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null

